$wb = $excel.ActiveWorkbook

# Sheet "2025"
$ws = $wb.Worksheets.Item("2025")
$ws.Range("N2").Value = 5744.096799946033
$ws.Range("O2").Value = 5589.635007435183

# Sheet "2030"
$ws = $wb.Worksheets.Item("2030")
$ws.Range("B2").Value = 4566.25257382453
$ws.Range("I2").Value = 35593.64721591155
$ws.Range("L2").Value = 53067.25365860503
$ws.Range("M2").Value = 18080.22084960085
$ws.Range("O2").Value = 9678.014132395467

# Sheet "2035"
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 2341.888254333185
$ws.Range("B2").Value = 6352.710371959407
$ws.Range("E2").Value = 53743.99346900621
$ws.Range("I2").Value = 47624.60274704682
$ws.Range("L2").Value = 53067.25365860503
$ws.Range("M2").Value = 21003.53454366937
$ws.Range("N2").Value = 12138.52420792558
$ws.Range("O2").Value = 11824.45056108007

# Sheet "2040"
$ws = $wb.Worksheets.Item("2040")
$ws.Range("A2").Value = 2341.888254333185
$ws.Range("B2").Value = 6352.710371959407
$ws.Range("E2").Value = 53743.99346900621
$ws.Range("I2").Value = 47624.60274704682
$ws.Range("L2").Value = 53067.25365860503
$ws.Range("M2").Value = 21003.53454366937
$ws.Range("N2").Value = 12225.049677232
$ws.Range("O2").Value = 11824.45056108007

# Sheet "2045"
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 5082.388487423812
$ws.Range("B2").Value = 6352.710371959407
$ws.Range("E2").Value = 53743.99346900621
$ws.Range("I2").Value = 47624.60274704682
$ws.Range("L2").Value = 53067.25365860503
$ws.Range("M2").Value = 21003.53454366937
$ws.Range("N2").Value = 12665.64083094451
$ws.Range("O2").Value = 13703.30132394448

# Sheet "2050"
$ws = $wb.Worksheets.Item("2050")
$ws.Range("A2").Value = 5082.388487423812
$ws.Range("B2").Value = 6352.710371959407
$ws.Range("E2").Value = 53743.99346900621
$ws.Range("I2").Value = 47624.60274704682
$ws.Range("L2").Value = 53067.25365860503
$ws.Range("M2").Value = 21003.53454366937
$ws.Range("N2").Value = 12665.64083094451
$ws.Range("O2").Value = 13703.30132394448
